# feat: add 2022-Q1 data
#
# - The sheet that used to be named "总计" (totals) is repurposed to hold the
#   new "2022-Q1" per-fund holdings table (same shape as the "2021-Q3" sheet).
# - A brand-new "总计" sheet is appended at the end, holding the running
#   totals table (one row per quarter), now with a "2022-Q1" row added above
#   the pre-existing "2021-Q3" row.

$wb = $excel.ActiveWorkbook

$q3Sheet = $wb.Worksheets.Item("2021-Q3")
$totalsSheet = $wb.Worksheets.Item("总计")

# A cell that already carries the bold/centered/bordered "header" style used
# by the "总计" table (and, after this edit, by the new per-fund table too) -
# we copy *formats only* from it so both new tables reuse the existing style
# instead of fabricating a new one. Do all the format-copying up front, while
# this donor cell is still untouched, before any values get overwritten.
$styleDonor = $totalsSheet.Range("B1")

# New "总计" sheet, appended after every existing sheet (i.e. after the old
# "总计" sheet, which is renamed to "2022-Q1" first below - sheet names must
# be unique, so the old sheet has to give up the "总计" name before the new
# sheet can claim it).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newTotals = $wb.Worksheets.Add($null, $lastSheet)

# ---------------------------------------------------------------------
# 1) Turn the old "总计" sheet into the new "2022-Q1" per-fund sheet.
# ---------------------------------------------------------------------
$totalsSheet.Name = "2022-Q1"
$q1Sheet = $totalsSheet

$newTotals.Name = "总计"

# Stamp the shared header style onto both destinations' header rows and index
# columns before writing any values into them.
$styleDonor.Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)
$styleDonor.Copy()
$q1Sheet.Range("A2:A11").PasteSpecial(-4122)

$styleDonor.Copy()
$newTotals.Range("B1:D1").PasteSpecial(-4122)
$styleDonor.Copy()
$newTotals.Range("A2:A3").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $q1Sheet.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

# row data: code, name, size, total position, position pct, held value, rank
$rows = @(
    @("012348", "天弘恒生科技指数型发起式证券投资基金（QDII）A", "38.10", "92.34", "7.52", "2.8651", 3),
    @("968029", "恒生指数基金M类人民币（对冲）份额",              "25.09", "97.94", "7.17", "1.7990", 4),
    @("012349", "天弘恒生科技指数型发起式证券投资基金（QDII）C", "14.77", "92.34", "7.52", "1.1107", 3),
    @("009562", "工银瑞信中国机会全球配置股票(QDII)美元",        "6.65",  "92.85", "3.00", "0.1995", 4),
    @("486001", "工银瑞信中国机会全球配置股票(QDII)",            "6.65",  "92.85", "3.00", "0.1995", 4),
    @("009563", "工银瑞信中国机会全球配置股票(QDII)港币",        "6.65",  "92.85", "3.00", "0.1995", 4),
    @("009225", "天弘中证中美互联网指数（QDII）A",                "1.84",  "94.90", "7.71", "0.1419", 5),
    @("002379", "工银瑞信香港中小盘股票（QDII）人民币",          "1.84",  "86.48", "3.53", "0.0650", 9),
    @("002380", "工银瑞信香港中小盘股票（QDII）美元",            "1.84",  "86.48", "3.53", "0.0650", 9),
    @("009226", "天弘中证中美互联网指数（QDII）C",                "0.59",  "94.90", "7.71", "0.0455", 5)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 2 + $i
    $row = $rows[$i]

    $q1Sheet.Cells.Item($r, 1).Value = $i
    # Fund codes are zero-padded numeric-looking strings (e.g. "012348") -
    # quote-prefix so the leading zero survives instead of being parsed away.
    $q1Sheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1Sheet.Cells.Item($r, 3).Value = $row[1]
    # Size / position / pct / value columns are stored as plain text in the
    # source data (e.g. "38.10"), not numbers - use a quote-prefix so Excel
    # doesn't silently coerce them to numeric cells.
    $q1Sheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1Sheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1Sheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $q1Sheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $q1Sheet.Cells.Item($r, 8).Value = $row[6]
}

$q1Sheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2) Fill in the new "总计" sheet with the running per-quarter totals.
# ---------------------------------------------------------------------
$newTotals.Cells.Item(1, 2).Value = "日期"
$newTotals.Cells.Item(1, 3).Value = "持有数量(只)"
$newTotals.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalsRows = @(
    @("2022-Q1", 10, 6.69),
    @("2021-Q3", 12, 4.1)
)

for ($i = 0; $i -lt $totalsRows.Count; $i++) {
    $r = 2 + $i
    $row = $totalsRows[$i]
    $newTotals.Cells.Item($r, 1).Value = $i
    $newTotals.Cells.Item($r, 2).Value = $row[0]
    $newTotals.Cells.Item($r, 3).Value = $row[1]
    $newTotals.Cells.Item($r, 4).Value = $row[2]
}

$newTotals.Range("A1").Select()
